# Remove the (now resolved / no longer needed) reviewer comment on slide 1:
#   "date and time should be in the same line?"
# This deletes ppt/comments/comment1.xml (and the now-dangling relationship
# to it + its Content-Types override), while leaving the comment author
# list (ppt/commentAuthors.xml) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $comment = $s.Comments.Item($i)
    if ($comment.Text -eq "date and time should be in the same line?") {
        $comment.Delete()
    }
}
